# Add the new weekly-ranking sheet '2025-12-31' after the last existing sheet
$wb = $excel.ActiveWorkbook
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = '2025-12-31'

# Bulk-write rank/title/volume values (columns A:C) as a single 2D array write
$values = New-Object 'object[,]' 101,3
$values[0,0] = 'rank'
$values[0,1] = 'title'
$values[0,2] = 'volume'
$values[1,0] = 1
$values[1,1] = '葬送のフリーレン'
$values[1,2] = 15
$values[2,0] = 2
$values[2,1] = '俺だけレベルアップな件'
$values[2,2] = 23
$values[3,0] = 3
$values[3,1] = '薬屋のひとりごと~猫猫の後宮謎解き手帳~'
$values[3,2] = 21
$values[4,0] = 4
$values[4,1] = 'ジョジョの奇妙な冒険 ザ・ジョジョランズ'
$values[4,2] = 7
$values[5,0] = 5
$values[5,1] = 'ダイヤモンドの功罪'
$values[5,2] = 9
$values[6,0] = 6
$values[6,1] = 'ONE PIECE'
$values[6,2] = 113
$values[7,0] = 7
$values[7,1] = 'ドンケツ第2章'
$values[7,2] = 16
$values[8,0] = 8
$values[8,1] = 'ウマ娘 シンデレラグレイ'
$values[8,2] = 22
$values[9,0] = 9
$values[9,1] = 'バトルスタディーズ'
$values[9,2] = 47
$values[10,0] = 10
$values[10,1] = 'ホタルの嫁入り'
$values[10,2] = 10
$values[11,0] = 11
$values[11,1] = '帝乃三姉妹は案外、チョロい。'
$values[11,2] = 17
$values[12,0] = 12
$values[12,1] = 'ネトオク男の楽しい異世界貿易'
$values[12,2] = 1
$values[13,0] = 13
$values[13,1] = 'ダーウィン事変'
$values[13,2] = 10
$values[14,0] = 14
$values[14,1] = '廻天のアルバス'
$values[14,2] = 7
$values[15,0] = 15
$values[15,1] = '規格外のダンジョン攻略者、実は異世界帰りの元勇者1'
$values[15,2] = 1
$values[16,0] = 16
$values[16,1] = 'THE BAND'
$values[16,2] = 3
$values[17,0] = 17
$values[17,1] = '変な家:'
$values[17,2] = 6
$values[18,0] = 18
$values[18,1] = '薬屋のひとりごと'
$values[18,2] = 16
$values[19,0] = 19
$values[19,1] = 'もう興味がないと離婚された令嬢の意外と楽しい新生活'
$values[19,2] = 5
$values[20,0] = 20
$values[20,1] = '災悪のアヴァロン~ゲーム最弱の悪役デブに転移したけど、俺だけ“やせれば強くてニューゲーム”な世界だったので、最速レベルアップ&破滅フラグ回避で影の英雄を目指します~'
$values[20,2] = 10
$values[21,0] = 21
$values[21,1] = '兼松先生、美味しゅうございますか?'
$values[21,2] = 1
$values[22,0] = 22
$values[22,1] = '呪われ公爵と捨てられた花嫁の最愛婚1'
$values[22,2] = 1
$values[23,0] = 23
$values[23,1] = '恋愛不感症―ホントはもっと感じたい―'
$values[23,2] = 4
$values[24,0] = 24
$values[24,1] = '修羅幼女の英雄譚~半端者と言われた傭兵、幼女に転生して成り上がる~1'
$values[24,2] = 1
$values[25,0] = 25
$values[25,1] = '無能と追放された最弱魔法剣士、呪いが解けたので最強へ成り上がる1'
$values[25,2] = 1
$values[26,0] = 26
$values[26,1] = '機動戦士ガンダム サンダーボルト'
$values[26,2] = 27
$values[27,0] = 27
$values[27,1] = '薫る花は凛と咲く'
$values[27,2] = 21
$values[28,0] = 28
$values[28,1] = '田舎者にはよくわかりません ぼんやり辺境伯令嬢は、断罪された公爵令息をお持ち帰りする'
$values[28,2] = 1
$values[29,0] = 29
$values[29,1] = 'わたくしに恋してください! ~ループ二回目の悪役令嬢ですが破滅回避のため誘惑します~'
$values[29,2] = 1
$values[30,0] = 30
$values[30,1] = '星廻りのレヴィア'
$values[30,2] = 1
$values[31,0] = 31
$values[31,1] = '今世では、ひとりで生きようと思います。そのはずが…1'
$values[31,2] = 1
$values[32,0] = 32
$values[32,1] = 'ファンタジーには馴染めない ~アラフォー男、ハードモード異世界に転移したけど結局無双~'
$values[32,2] = 1
$values[33,0] = 33
$values[33,1] = '暴君王の初恋'
$values[33,2] = 3
$values[34,0] = 34
$values[34,1] = '恋だ獣'
$values[34,2] = 3
$values[35,0] = 35
$values[35,1] = '僕の彼女はデッカワイイ'
$values[35,2] = 1
$values[36,0] = 36
$values[36,1] = '田舎の黒ギャルJKと結婚しました'
$values[36,2] = 1
$values[37,0] = 37
$values[37,1] = '同居している剣聖の女師匠が可愛すぎて毎日幸せです'
$values[37,2] = 1
$values[38,0] = 38
$values[38,1] = '冒険者絶対殺すダンジョン'
$values[38,2] = 1
$values[39,0] = 39
$values[39,1] = '執事ですがなにか?~幼馴染のパワハラ皇女と絶縁したら、隣国の向日葵王女に拾われたのでこの身を捧げます~1'
$values[39,2] = 1
$values[40,0] = 40
$values[40,1] = '高嶺のハナさん'
$values[40,2] = 13
$values[41,0] = 41
$values[41,1] = 'ザ・ファブル The third secret'
$values[41,2] = 3
$values[42,0] = 42
$values[42,1] = 'みいちゃんと山田さん'
$values[42,2] = 5
$values[43,0] = 43
$values[43,1] = 'ケンガンオメガ'
$values[43,2] = 31
$values[44,0] = 44
$values[44,1] = 'ブラックナイトパレード'
$values[44,2] = 11
$values[45,0] = 45
$values[45,1] = '最弱テイマーはゴミ拾いの旅を始めました。@COMIC'
$values[45,2] = 8
$values[46,0] = 46
$values[46,1] = '青の祓魔師'
$values[46,2] = 33
$values[47,0] = 47
$values[47,1] = '龍と苺'
$values[47,2] = 23
$values[48,0] = 48
$values[48,1] = 'レッドブルー'
$values[48,2] = 16
$values[49,0] = 49
$values[49,1] = 'ワンパンマン'
$values[49,2] = 35
$values[50,0] = 50
$values[50,1] = 'SAKAMOTO DAYS'
$values[50,2] = 25
$values[51,0] = 51
$values[51,1] = 'アオのハコ'
$values[51,2] = 23
$values[52,0] = 52
$values[52,1] = '魔女の婚姻'
$values[52,2] = 1
$values[53,0] = 53
$values[53,1] = '成り代わり令嬢のループライン'
$values[53,2] = 1
$values[54,0] = 54
$values[54,1] = '初恋の少年は冷徹騎士に豹変していました'
$values[54,2] = 1
$values[55,0] = 55
$values[55,1] = 'いぬみみ'
$values[55,2] = 1
$values[56,0] = 56
$values[56,1] = 'フェアリーテイル・クロニクル ~空気読まない異世界ライフ~'
$values[56,2] = 1
$values[57,0] = 57
$values[57,1] = '大正シンデレラ~秘密の恋は髪が伸びるまで~'
$values[57,2] = 2
$values[58,0] = 58
$values[58,1] = '離縁は致しかねます!'
$values[58,2] = 6
$values[59,0] = 59
$values[59,1] = '僕の彼女はデッカワイイ'
$values[59,2] = 2
$values[60,0] = 60
$values[60,1] = '僕の彼女はデッカワイイ'
$values[60,2] = 3
$values[61,0] = 61
$values[61,1] = '日々は過ぎれど飯うまし'
$values[61,2] = 1
$values[62,0] = 62
$values[62,1] = '朧の花嫁~かりそめの婚約は、青く、甘く~'
$values[62,2] = 1
$values[63,0] = 63
$values[63,1] = '魔法歌姫マジカルギンガ 第25話'
$values[63,2] = 25
$values[64,0] = 64
$values[64,1] = 'ブルーピリオド'
$values[64,2] = 18
$values[65,0] = 65
$values[65,1] = '焼いてるふたり'
$values[65,2] = 22
$values[66,0] = 66
$values[66,1] = '怨み屋本舗DIABLO'
$values[66,2] = 13
$values[67,0] = 67
$values[67,1] = '神様のバレー'
$values[67,2] = 39
$values[68,0] = 68
$values[68,1] = '盤上のオリオン'
$values[68,2] = 8
$values[69,0] = 69
$values[69,1] = 'WIND BREAKER'
$values[69,2] = 24
$values[70,0] = 70
$values[70,1] = 'ブルーロック'
$values[70,2] = 36
$values[71,0] = 71
$values[71,1] = '桃源暗鬼'
$values[71,2] = 27
$values[72,0] = 72
$values[72,1] = '魔入りました!入間くん'
$values[72,2] = 46
$values[73,0] = 73
$values[73,1] = 'トニカクカワイイ'
$values[73,2] = 34
$values[74,0] = 74
$values[74,1] = 'たわら猫とまちがい人生'
$values[74,2] = 1
$values[75,0] = 75
$values[75,1] = '転生したら小魚だったけど龍になれるらしいので頑張ります'
$values[75,2] = 1
$values[76,0] = 76
$values[76,1] = '恋検'
$values[76,2] = 1
$values[77,0] = 77
$values[77,1] = 'ブラック嫁によろしく!'
$values[77,2] = 1
$values[78,0] = 78
$values[78,1] = '初恋相手の兄に嫁ぎました'
$values[78,2] = 6
$values[79,0] = 79
$values[79,1] = '田舎の黒ギャルJKと結婚しました'
$values[79,2] = 2
$values[80,0] = 80
$values[80,1] = '田舎の黒ギャルJKと結婚しました'
$values[80,2] = 3
$values[81,0] = 81
$values[81,1] = '同居している剣聖の女師匠が可愛すぎて毎日幸せです'
$values[81,2] = 2
$values[82,0] = 82
$values[82,1] = '同居している剣聖の女師匠が可愛すぎて毎日幸せです'
$values[82,2] = 3
$values[83,0] = 83
$values[83,1] = '冒険者絶対殺すダンジョン'
$values[83,2] = 2
$values[84,0] = 84
$values[84,1] = '冒険者絶対殺すダンジョン'
$values[84,2] = 3
$values[85,0] = 85
$values[85,1] = '東方Project二次創作シリーズ 紅魔館の女たち'
$values[85,2] = 1
$values[86,0] = 86
$values[86,1] = 'おっさんはうぜぇぇぇんだよ!ってギルドから追放したくせに、後から復帰要請を出されても遅い。最高の仲間と出会った俺はこっちで最強を目指す! コミック版'
$values[86,2] = 4
$values[87,0] = 87
$values[87,1] = 'ちひろさん'
$values[87,2] = 10
$values[88,0] = 88
$values[88,1] = '凍牌 コールドガール'
$values[88,2] = 11
$values[89,0] = 89
$values[89,1] = '実は俺、最強でした?'
$values[89,2] = 18
$values[90,0] = 90
$values[90,1] = '失格紋の最強賢者 ~世界最強の賢者が更に強くなるために転生しました~'
$values[90,2] = 33
$values[91,0] = 91
$values[91,1] = '勇者パーティを追い出された器用貧乏 ~パーティ事情で付与術士をやっていた剣士、万能へと至る~'
$values[91,2] = 17
$values[92,0] = 92
$values[92,1] = 'めしぬま。'
$values[92,2] = 15
$values[93,0] = 93
$values[93,1] = 'MIX'
$values[93,2] = 24
$values[94,0] = 94
$values[94,1] = '陸奥圓明流異界伝 修羅の紋 ムツさんはチョー強い?!'
$values[94,2] = 16
$values[95,0] = 95
$values[95,1] = '青のミブロー新選組編ー'
$values[95,2] = 8
$values[96,0] = 96
$values[96,1] = 'ワールドトリガー'
$values[96,2] = 29
$values[97,0] = 97
$values[97,1] = '離婚予定の契約婚なのに、冷酷公爵様に執着されています'
$values[97,2] = 21
$values[98,0] = 98
$values[98,1] = '田舎者にはよくわかりません ぼんやり辺境伯令嬢は、断罪された公爵令息をお持ち帰りする'
$values[98,2] = 2
$values[99,0] = 99
$values[99,1] = '田舎者にはよくわかりません ぼんやり辺境伯令嬢は、断罪された公爵令息をお持ち帰りする'
$values[99,2] = 3
$values[100,0] = 100
$values[100,1] = '魔女の婚姻'
$values[100,2] = 2

$ws.Range('A1:C101').Value2 = $values

# Copy the bold/centered/bordered header style (xf s="1") from the last existing sheet
$lastSheet.Range('A1:D1').Copy()
$ws.Range('A1:D1').PasteSpecial(-4122)

# Copy the highlighted-volume fill style (xf s="2") onto the flagged volume cells
$fillSourceCell = $lastSheet.Range('C17')
$fillSourceCell.Copy()
$highlightRows = @(13,16,17,22,23,25,26,29,30,31,32,33,34,35,36,37,38,39,40,42,53,54,55,56,57,58,60,61,62,63,75,76,77,78,80,81,82,83,84,85,86,99,100,101)
foreach ($r in $highlightRows) {
  $ws.Range('C' + $r).PasteSpecial(-4122)
}

$excel.CutCopyMode = $false

$ws.Range('A1').Select()
